# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1079
    5  = 3073
    7  = 2389
    9  = 117
    10 = 1
    11 = 1187
    15 = 1070
    16 = 287
    17 = 307
    20 = 101
    21 = 59
    23 = 5
    24 = 31
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
